# Update column F ("dSF") values for specific rows, per repull of data /
# recalculation of mean (see commit message: "repull data, push all data,
# mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 4
    7  = 4
    9  = -9
    10 = -3
    11 = -1
    12 = 9
    16 = 2
    19 = -2
    20 = 4
    21 = -3
    22 = -7
    23 = -3
    25 = -2
    26 = -1
    28 = -3
    31 = 3
    32 = -7
    36 = -1
    37 = -4
    38 = 4
    52 = -3
    53 = 6
    56 = -2
    57 = 6
    59 = -2
    60 = 2
    61 = 1
    62 = 2
    65 = 1
    67 = 3
    70 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
